$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data range is treated as text so values like "1.00" or "588.49" are not
# auto-converted to numbers, matching the original inline-string cell formatting.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "63.552.43"
$ws.Cells.Item(2, 5).Value = "  -1.26%  "
$ws.Cells.Item(3, 4).Value = "3.071.28"
$ws.Cells.Item(3, 5).Value = "  -3.35%  "
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  -0.06%  "
$ws.Cells.Item(5, 4).Value = "588.49"
$ws.Cells.Item(5, 5).Value = "  -0.97%  "
$ws.Cells.Item(6, 4).Value = "155.31"
$ws.Cells.Item(6, 5).Value = "  +4.42%  "
$ws.Cells.Item(7, 5).Value = "  -0.07%  "
$ws.Cells.Item(8, 5).Value = "  +0.39%  "
$ws.Cells.Item(9, 4).Value = "3.068.56"
$ws.Cells.Item(9, 5).Value = "  -3.09%  "
$ws.Cells.Item(10, 4).Value = "0.156"
$ws.Cells.Item(10, 5).Value = "  -3.98%  "
$ws.Cells.Item(11, 5).Value = "  -1.57%  "
$ws.Cells.Item(12, 4).Value = "0.449"
$ws.Cells.Item(12, 5).Value = "  -2.85%  "
$ws.Cells.Item(13, 4).Value = "36.84"
$ws.Cells.Item(13, 5).Value = "  -2.45%  "
$ws.Cells.Item(14, 5).Value = "  -4.79%  "
$ws.Cells.Item(15, 2).Value = "TRON"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(15, 4).Value = "0.119"
$ws.Cells.Item(15, 5).Value = "  -2.21%  "
$ws.Cells.Item(16, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(16, 4).Value = "3.577.06"
$ws.Cells.Item(16, 5).Value = "  -3.33%  "
$ws.Cells.Item(17, 4).Value = "63.572.99"
$ws.Cells.Item(17, 5).Value = "  -0.96%  "
$ws.Cells.Item(18, 4).Value = "7.12"
$ws.Cells.Item(18, 5).Value = "  -2.85%  "
$ws.Cells.Item(19, 4).Value = "3.072.69"
$ws.Cells.Item(19, 5).Value = "  -3.09%  "
$ws.Cells.Item(20, 4).Value = "469.30"
$ws.Cells.Item(20, 5).Value = "  -0.33%  "
$ws.Cells.Item(21, 4).Value = "14.29"
$ws.Cells.Item(21, 5).Value = "  -1.92%  "
$ws.Cells.Item(22, 4).Value = "0.703"
$ws.Cells.Item(22, 5).Value = "  -4.86%  "
$ws.Cells.Item(23, 4).Value = "7.49"
$ws.Cells.Item(23, 5).Value = "  -2.87%  "
$ws.Cells.Item(24, 5).Value = "  -1.73%  "
$ws.Cells.Item(25, 4).Value = "80.50"
$ws.Cells.Item(25, 5).Value = "  -1.37%  "
$ws.Cells.Item(26, 4).Value = "12.77"
$ws.Cells.Item(26, 5).Value = "  -3.59%  "
$ws.Cells.Item(27, 4).Value = "10.36"
$ws.Cells.Item(27, 5).Value = "  +2.35%  "
$ws.Cells.Item(28, 4).Value = "0.998"
$ws.Cells.Item(28, 5).Value = "  -0.21%  "
$ws.Cells.Item(29, 4).Value = "7.41"
$ws.Cells.Item(29, 5).Value = "  +1.76%  "
$ws.Cells.Item(30, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(30, 4).Value = "1.00"
$ws.Cells.Item(30, 5).Value = "  +0.00%  "
$ws.Cells.Item(31, 2).Value = "PancakeSwap"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(31, 4).Value = "2.65"
$ws.Cells.Item(31, 5).Value = "  -3.20%  "
$ws.Cells.Item(32, 5).Value = "  -5.52%  "
$ws.Cells.Item(33, 4).Value = "27.10"
$ws.Cells.Item(33, 5).Value = "  -4.75%  "
$ws.Cells.Item(34, 5).Value = "  -5.10%  "
$ws.Cells.Item(35, 4).Value = "0.0₃0818"
$ws.Cells.Item(35, 5).Value = "  -5.20%  "
$ws.Cells.Item(36, 5).Value = "  -2.46%  "
$ws.Cells.Item(37, 4).Value = "5.97"
$ws.Cells.Item(37, 5).Value = "  -4.14%  "
$ws.Cells.Item(38, 4).Value = "3.26"
$ws.Cells.Item(38, 5).Value = "  -2.21%  "
$ws.Cells.Item(39, 4).Value = "2.20"
$ws.Cells.Item(39, 5).Value = "  -5.57%  "
$ws.Cells.Item(40, 4).Value = "50.57"
$ws.Cells.Item(40, 5).Value = "  -1.90%  "
$ws.Cells.Item(41, 4).Value = "9.18"
$ws.Cells.Item(41, 5).Value = "  -1.76%  "
$ws.Cells.Item(42, 4).Value = "435.93"
$ws.Cells.Item(42, 5).Value = "  -7.48%  "
$ws.Cells.Item(43, 5).Value = "  -3.85%  "
$ws.Cells.Item(44, 4).Value = "40.44"
$ws.Cells.Item(44, 5).Value = "  +1.72%  "
$ws.Cells.Item(45, 5).Value = "  +2.67%  "
$ws.Cells.Item(46, 4).Value = "0.0358"
$ws.Cells.Item(46, 5).Value = "  -4.91%  "
$ws.Cells.Item(47, 4).Value = "2.794.54"
$ws.Cells.Item(47, 5).Value = "  -4.16%  "
$ws.Cells.Item(48, 4).Value = "128.98"
$ws.Cells.Item(48, 5).Value = "  -3.12%  "
$ws.Cells.Item(50, 4).Value = "25.01"
$ws.Cells.Item(50, 5).Value = "  +2.03%  "
$ws.Cells.Item(51, 4).Value = "2.21"
$ws.Cells.Item(51, 5).Value = "  -2.24%  "

# Restore default (Normal) style so no stray number-format styling is left on cells,
# keeping the style/formatting identical to the original workbook.
$dataRange.Style = "Normal"

